$wb = $excel.ActiveWorkbook

# Sheets that contain the event data table: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 7735
    $ws.Range("G3").Value = 75

    $ws.Range("F9").Value = 5975

    $ws.Range("F12").Value = 30

    $ws.Range("F13").Value = 1817

    $ws.Range("F14").Value = 1334

    $ws.Range("F16").Value = 629

    $ws.Range("F17").Value = 145
}
